$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these cells keep their original text representation (matching the
# source inlineStr cells) rather than being auto-converted to numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.779.97"
$ws.Range("E2").Value = "  -1.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.635.60"
$ws.Range("E3").Value = "  -1.85%  "
$ws.Range("E4").Value = "  +0.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.78"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3861"
$ws.Range("E7").Value = "  -1.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3799"
$ws.Range("E8").Value = "  -2.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.42"
$ws.Range("E9").Value = "  -2.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.316"
$ws.Range("E10").Value = "  -4.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.003"
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08348"
$ws.Range("E12").Value = "  -2.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.57"
$ws.Range("E13").Value = "  -2.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.931"
$ws.Range("E14").Value = "  -4.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.756"
$ws.Range("E15").Value = "  -4.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001299"
$ws.Range("E16").Value = "  -1.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.636.63"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.26"
$ws.Range("E18").Value = "  -2.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06935"
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.28"
$ws.Range("E20").Value = "  -4.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.830"
$ws.Range("E21").Value = "  -2.69%  "
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.44"
$ws.Range("E23").Value = "  -2.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.781.17"
$ws.Range("E24").Value = "  -1.39%  "
$ws.Range("E25").Value = "  -2.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.844"
$ws.Range("E26").Value = "  -11.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.70"
$ws.Range("E27").Value = "  -3.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.02"
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.462"
$ws.Range("E29").Value = "  +3.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.00"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.758"
$ws.Range("E31").Value = "  -1.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.487"
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.817.87"
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07923"
$ws.Range("E34").Value = "  -2.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9765"
$ws.Range("E35").Value = "  -7.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02874"
$ws.Range("E36").Value = "  -5.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.533"
$ws.Range("E37").Value = "  -3.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2639"
$ws.Range("E38").Value = "  -3.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.34"
$ws.Range("E39").Value = "  -9.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09059"
$ws.Range("E40").Value = "  -1.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7421"
$ws.Range("E41").Value = "  -2.81%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.414"
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.16"
$ws.Range("E43").Value = "  -3.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.50"
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6831"
$ws.Range("E45").Value = "  -3.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.390"
$ws.Range("E46").Value = "  -5.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.061"
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.002"
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08194"
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.55"
$ws.Range("E50").Value = "  -1.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.210"
$ws.Range("E51").Value = "  -3.08%  "
